$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrC = New-Object 'object[,]' 24,1
$arrC[0,0] = 5.557378622606144
$arrC[1,0] = 5.565838940101819
$arrC[2,0] = 5.571260976727778
$arrC[3,0] = 5.573528031764074
$arrC[4,0] = 5.573907959801817
$arrC[5,0] = 5.571291317640463
$arrC[6,0] = 5.560248764541003
$arrC[7,0] = 5.540382306902557
$arrC[8,0] = 5.526854165953184
$arrC[9,0] = 5.520927028655699
$arrC[10,0] = 5.51871484225228
$arrC[11,0] = 5.519189845239573
$arrC[12,0] = 5.520744385468318
$arrC[13,0] = 5.521700781770657
$arrC[14,0] = 5.527246055146943
$arrC[15,0] = 5.530705771914693
$arrC[16,0] = 5.532717084351817
$arrC[17,0] = 5.533401761557568
$arrC[18,0] = 5.530335269373246
$arrC[19,0] = 5.520286905562864
$arrC[20,0] = 5.513907788930098
$arrC[21,0] = 5.517295344134813
$arrC[22,0] = 5.530502704164615
$arrC[23,0] = 5.545567613402174
$ws.Range("C2:C25").Value = $arrC

$arrD = New-Object 'object[,]' 24,1
$arrD[0,0] = 6.406566698564359
$arrD[1,0] = 6.390393480460238
$arrD[2,0] = 6.381128704289981
$arrD[3,0] = 6.377522936581746
$arrD[4,0] = 6.376934517024516
$arrD[5,0] = 6.381079385518704
$arrD[6,0] = 6.40085344121947
$arrD[7,0] = 6.444806252535138
$arrD[8,0] = 6.480107067672117
$arrD[9,0] = 6.496786397332429
$arrD[10,0] = 6.503188714144716
$arrD[11,0] = 6.501806075711157
$arrD[12,0] = 6.497311409679806
$arrD[13,0] = 6.494569432783695
$arrD[14,0] = 6.479029267363673
$arrD[15,0] = 6.469652737761003
$arrD[16,0] = 6.464318228164139
$arrD[17,0] = 6.462522214976907
$arrD[18,0] = 6.470644840025725
$arrD[19,0] = 6.498629286920865
$arrD[20,0] = 6.517419801761447
$arrD[21,0] = 6.507346134160742
$arrD[22,0] = 6.470196135366536
$arrD[23,0] = 6.432376467434326
$ws.Range("D2:D25").Value = $arrD

$arrE = New-Object 'object[,]' 24,1
$arrE[0,0] = 12.53289745328044
$arrE[1,0] = 12.49979517383337
$arrE[2,0] = 12.48208417343762
$arrE[3,0] = 12.47552810870916
$arrE[4,0] = 12.47447952570199
$arrE[5,0] = 12.48199307358454
$arrE[6,0] = 12.52094386103263
$arrE[7,0] = 12.61786349930744
$arrE[8,0] = 12.70128051368273
$arrE[9,0] = 12.74180194207017
$arrE[10,0] = 12.7575089974247
$arrE[11,0] = 12.75411021029441
$arrE[12,0] = 12.74308695286921
$arrE[13,0] = 12.73638185428245
$arrE[14,0] = 12.6986835002244
$arrE[15,0] = 12.67621085636371
$arrE[16,0] = 12.66352802036852
$arrE[17,0] = 12.6592757644055
$arrE[18,0] = 12.67857802969121
$arrE[19,0] = 12.74631497793738
$arrE[20,0] = 12.79269377049733
$arrE[21,0] = 12.76775028222931
$arrE[22,0] = 12.67750709188873
$arrE[23,0] = 12.58947506639407
$ws.Range("E2:E25").Value = $arrE

$arrF = New-Object 'object[,]' 24,1
$arrF[0,0] = 46.378627352927
$arrF[1,0] = 46.27821239909783
$arrF[2,0] = 46.23045381635081
$arrF[3,0] = 46.21448703429792
$arrF[4,0] = 46.21204681257926
$arrF[5,0] = 46.23022433285585
$arrF[6,0] = 46.34112167967328
$arrF[7,0] = 46.66864731102498
$arrF[8,0] = 46.97584886849857
$arrF[9,0] = 47.12986406255506
$arrF[10,0] = 47.19021463260532
$arrF[11,0] = 47.17712725636662
$arrF[12,0] = 47.13478863571384
$arrF[13,0] = 47.1091184047033
$arrF[14,0] = 46.96606880779262
$arrF[15,0] = 46.88195141080875
$arrF[16,0] = 46.8349138877943
$arrF[17,0] = 46.81921937208506
$arrF[18,0] = 46.89076684954956
$arrF[19,0] = 47.14716967004427
$arrF[20,0] = 47.32655315747776
$arrF[21,0] = 47.22974090588685
$arrF[22,0] = 46.88677726809077
$arrF[23,0] = 46.56829914222566
$ws.Range("F2:F25").Value = $arrF

$arrG = New-Object 'object[,]' 24,1
$arrG[0,0] = 3.72147277115217
$arrG[1,0] = 3.725632090746696
$arrG[2,0] = 3.728316095681631
$arrG[3,0] = 3.729442713293437
$arrG[4,0] = 3.729631776069716
$arrG[5,0] = 3.728331156403552
$arrG[6,0] = 3.72287996860076
$arrG[7,0] = 3.713216962623422
$arrG[8,0] = 3.706735012936941
$arrG[9,0] = 3.703918466402405
$arrG[10,0] = 3.702870772766887
$arrG[11,0] = 3.703095575061048
$arrG[12,0] = 3.703831894537529
$arrG[13,0] = 3.704285365298201
$arrG[14,0] = 3.706921728620648
$arrG[15,0] = 3.70857280132004
$arrG[16,0] = 3.709534898525218
$arrG[17,0] = 3.709862789253433
$arrG[18,0] = 3.708395754849435
$arrG[19,0] = 3.703615108511926
$arrG[20,0] = 3.700600617587044
$arrG[21,0] = 3.702199491023487
$arrG[22,0] = 3.708475757451046
$arrG[23,0] = 3.715722017962313
$ws.Range("G2:G25").Value = $arrG

$arrI = New-Object 'object[,]' 24,1
$arrI[0,0] = 32.71530555892914
$arrI[1,0] = 32.67285520587767
$arrI[2,0] = 32.65642392917027
$arrI[3,0] = 32.65214657785641
$arrI[4,0] = 32.65158224242558
$arrI[5,0] = 32.65635645782532
$arrI[6,0] = 32.69866760233944
$arrI[7,0] = 32.8581504670208
$arrI[8,0] = 33.0219241168471
$arrI[9,0] = 33.1064904038155
$arrI[10,0] = 33.13995189821777
$arrI[11,0] = 33.13268157960858
$arrI[12,0] = 33.10921454109072
$arrI[13,0] = 33.0950272685192
$arrI[14,0] = 33.01659930544569
$arrI[15,0] = 32.9710583113356
$arrI[16,0] = 32.94581238937441
$arrI[17,0] = 32.93742762579942
$arrI[18,0] = 32.97580816746071
$arrI[19,0] = 33.11606843530276
$arrI[20,0] = 33.21611292368517
$arrI[21,0] = 33.16195451601925
$arrI[22,0] = 32.97365784077356
$arrI[23,0] = 32.80681001458179
$ws.Range("I2:I25").Value = $arrI

$arrJ = New-Object 'object[,]' 24,1
$arrJ[0,0] = 10.41249716900742
$arrJ[1,0] = 10.42335079099023
$arrJ[2,0] = 10.43191938287488
$arrJ[3,0] = 10.43588928883696
$arrJ[4,0] = 10.43657734275408
$arrJ[5,0] = 10.43197098754911
$arrJ[6,0] = 10.41584381669688
$arrJ[7,0] = 10.39935914015863
$arrJ[8,0] = 10.39651367207409
$arrJ[9,0] = 10.39723631244238
$arrJ[10,0] = 10.39780016024804
$arrJ[11,0] = 10.39766581794823
$arrJ[12,0] = 10.39727688415353
$arrJ[13,0] = 10.39707644555789
$arrJ[14,0] = 10.39650705013126
$arrJ[15,0] = 10.39667451185479
$arrJ[16,0] = 10.3969606864373
$arrJ[17,0] = 10.39709018116651
$arrJ[18,0] = 10.39663703447463
$arrJ[19,0] = 10.39738324701056
$arrJ[20,0] = 10.39956241518341
$arrJ[21,0] = 10.39824456903529
$arrJ[22,0] = 10.39665338651925
$arrJ[23,0] = 10.40219357564587
$ws.Range("J2:J25").Value = $arrJ

$arrK = New-Object 'object[,]' 24,1
$arrK[0,0] = 22.27142012897855
$arrK[1,0] = 21.88596838644426
$arrK[2,0] = 21.65196084934779
$arrK[3,0] = 21.55739920333044
$arrK[4,0] = 21.54174929527755
$arrK[5,0] = 21.65068215535835
$arrK[6,0] = 22.13804147590085
$arrK[7,0] = 23.10921483083354
$arrK[8,0] = 23.82470589129908
$arrK[9,0] = 24.14910809572487
$arrK[10,0] = 24.27167465257711
$arrK[11,0] = 24.24529198752564
$arrK[12,0] = 24.15919796459271
$arrK[13,0] = 24.10642320423225
$arrK[14,0] = 23.80347315949011
$arrK[15,0] = 23.61725818666944
$arrK[16,0] = 23.51005778469931
$arrK[17,0] = 23.47374916047849
$arrK[18,0] = 23.6370917699302
$arrK[19,0] = 24.18449434749468
$arrK[20,0] = 24.54058303913626
$arrK[21,0] = 24.35072372039875
$arrK[22,0] = 23.62812544279568
$arrK[23,0] = 22.84561128688647
$ws.Range("K2:K25").Value = $arrK

$arrN = New-Object 'object[,]' 24,1
$arrN[0,0] = 19.37751929838118
$arrN[1,0] = 19.45890956508108
$arrN[2,0] = 19.51094575903455
$arrN[3,0] = 19.53267172461789
$arrN[4,0] = 19.53631082452302
$arrN[5,0] = 19.51123665146876
$arrN[6,0] = 19.40515599326881
$arrN[7,0] = 19.21339148281083
$arrN[8,0] = 19.08226626488034
$arrN[9,0] = 19.02470227580847
$arrN[10,0] = 19.0032017676841
$arrN[11,0] = 19.00781907994132
$arrN[12,0] = 19.02292746165022
$arrN[13,0] = 19.03222048985356
$arrN[14,0] = 19.08606998750145
$arrN[15,0] = 19.11963751909108
$arrN[16,0] = 19.1391411052513
$arrN[17,0] = 19.14577848753878
$arrN[18,0] = 19.11604388404995
$arrN[19,0] = 19.01848170083081
$arrN[20,0] = 18.95645353717795
$arrN[21,0] = 18.98940115846542
$arrN[22,0] = 19.11766792724448
$arrN[23,0] = 19.26354348996357
$ws.Range("N2:N25").Value = $arrN

